$d = $word.ActiveDocument

# Locate the heading text that needs to change.
$rng = $d.Content
$found = $rng.Find.Execute("What is wrong with the following code", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $paraStart = $rng.Start

    # Replace the whole matched phrase with the new wording.
    $rng.Text = "What security issue is prevalent in the code below"
    $newEnd = $rng.End

    # Force a run-break between "What " and "security issue..." by
    # dropping a throwaway bookmark at that boundary and then removing it -
    # the text either side of a bookmark is always kept in separate runs.
    $splitPoint = $paraStart + 5
    $splitRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("__TempSplit", $splitRange)

    # Re-create the "_GoBack" bookmark right after the new sentence (before
    # " and how would you fix it?"). Because "_GoBack" already exists lower
    # down in the document (after the PHP snippet), adding it here moves it
    # from its old location to this new one.
    $goBackRange = $d.Range($newEnd, $newEnd)
    $d.Bookmarks.Add("_GoBack", $goBackRange)

    # Drop the temporary bookmark now that the run split is baked in.
    $d.Bookmarks.Item("__TempSplit").Delete()
}
